# The workbook's species-observation rows 3 and 4 got reordered/swapped.
# Row 3 must take on what used to be row 4's record, and vice versa.
# All other rows are untouched.
#
# Note: a handful of text cells hold purely-numeric or date-shaped text
# (e.g. "1", "2013-05-16"). Assigning such text straight to .Value2 makes
# Excel auto-convert it to a real number/date. To keep them as genuine
# text (matching the source data), those are entered with a leading
# apostrophe (forces text entry, like typing it in the UI) and the cell
# style is then reset to "Normal" so no quote-prefix formatting lingers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 becomes the old "Robust tickgnagare / Dorcatoma robusta" record ---
$ws.Range("A3").Value2 = 112141528
$ws.Range("B3").Value2 = 4755
$ws.Range("D3").Value2 = "LC"
$ws.Range("E3").Value2 = 100857
$ws.Range("F3").Value2 = "Robust tickgnagare"
$ws.Range("G3").Value2 = "Dorcatoma robusta"
$ws.Range("H3").Value2 = "Strand, 1938"

$ws.Range("I3").Value2 = "'1"
$ws.Range("I3").Style = "Normal"

$ws.Range("J3").Value2 = "ex."
$ws.Range("K3").Value2 = "imago/adult"
$ws.Range("M3").Value2 = ""
$ws.Range("N3").Value2 = "fönsterfälla"
$ws.Range("Q3").Value2 = 445824
$ws.Range("R3").Value2 = 6205171

$ws.Range("Y3").Value2 = "'2013-05-16"
$ws.Range("Y3").Style = "Normal"

$ws.Range("AA3").Value2 = "'2013-05-24"
$ws.Range("AA3").Style = "Normal"

$ws.Range("AI3").Value2 = "i gles tallskog"
$ws.Range("AO3").Value2 = "på nydöd tall, delvis barklös"
$ws.Range("AQ3").Value2 = "Nils Otto Nilsson"
$ws.Range("AR3").Value2 = "NON 04741"

# --- Row 4 becomes the old "Åkerväddsantennmal / Nemophora metallica" record ---
$ws.Range("A4").Value2 = 112156959
$ws.Range("B4").Value2 = 39455
$ws.Range("D4").Value2 = "NT"
$ws.Range("E4").Value2 = 102471
$ws.Range("F4").Value2 = "Åkerväddsantennmal"
$ws.Range("G4").Value2 = "Nemophora metallica"
$ws.Range("H4").Value2 = "(Poda, 1761)"

$ws.Range("I4").Value2 = "'"
$ws.Range("I4").Style = "Normal"

$ws.Range("J4").Value2 = ""
$ws.Range("K4").Value2 = ""
$ws.Range("M4").Value2 = "födosökande"
$ws.Range("N4").Value2 = ""
$ws.Range("Q4").Value2 = 445825
$ws.Range("R4").Value2 = 6205212

$ws.Range("Y4").Value2 = "'2013-07-11"
$ws.Range("Y4").Style = "Normal"

$ws.Range("AA4").Value2 = "'2013-07-11"
$ws.Range("AA4").Style = "Normal"

$ws.Range("AI4").Value2 = "i tallskogsbryn"
$ws.Range("AO4").Value2 = "på blmr av åkervädd"
$ws.Range("AQ4").Value2 = ""
$ws.Range("AR4").Value2 = ""
